{"js": "// Replace each three-digit-by-one-digit multiplication equation in the\n// document's table cells with its new value, per the commit's regenerated\n// answer set. Every equation string is unique in the document, so a plain\n// exact-text search/replace per pair is unambiguous.\nconst replacements = [\n  [\"611\u00d78=4888\", \"258\u00d76=1548\"],\n  [\"539\u00d72=1078\", \"247\u00d72=494\"],\n  [\"212\u00d78=1696\", \"363\u00d78=2904\"],\n  [\"711\u00d76=4266\", \"974\u00d79=8766\"],\n  [\"818\u00d78=6544\", \"693\u00d78=5544\"],\n  [\"470\u00d78=3760\", \"587\u00d73=1761\"],\n  [\"947\u00d75=4735\", \"977\u00d75=4885\"],\n  [\"715\u00d73=2145\", \"599\u00d78=4792\"],\n  [\"725\u00d76=4350\", \"216\u00d79=1944\"],\n  [\"692\u00d74=2768\", \"902\u00d76=5412\"],\n  [\"305\u00d72=610\", \"824\u00d77=5768\"],\n  [\"443\u00d79=3987\", \"460\u00d79=4140\"],\n  [\"420\u00d75=2100\", \"833\u00d73=2499\"],\n  [\"885\u00d79=7965\", \"311\u00d73=933\"],\n  [\"984\u00d72=1968\", \"803\u00d78=6424\"],\n  [\"529\u00d77=3703\", \"912\u00d79=8208\"],\n  [\"340\u00d72=680\", \"223\u00d79=2007\"],\n  [\"309\u00d74=1236\", \"987\u00d73=2961\"],\n  [\"107\u00d72=214\", \"318\u00d73=954\"],\n  [\"436\u00d77=3052\", \"799\u00d74=3196\"],\n  [\"768\u00d77=5376\", \"293\u00d75=1465\"],\n  [\"628\u00d74=2512\", \"690\u00d73=2070\"],\n  [\"852\u00d76=5112\", \"748\u00d76=4488\"],\n  [\"132\u00d79=1188\", \"956\u00d76=5736\"],\n  [\"185\u00d75=925\", \"350\u00d73=1050\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication equation in the\n# document's table cells with its new value, per the commit's regenerated\n# answer set. Every equation string is unique in the document, so a plain\n# exact-text Find/Replace per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"611\u00d78=4888\", \"258\u00d76=1548\"),\n    @(\"539\u00d72=1078\", \"247\u00d72=494\"),\n    @(\"212\u00d78=1696\", \"363\u00d78=2904\"),\n    @(\"711\u00d76=4266\", \"974\u00d79=8766\"),\n    @(\"818\u00d78=6544\", \"693\u00d78=5544\"),\n    @(\"470\u00d78=3760\", \"587\u00d73=1761\"),\n    @(\"947\u00d75=4735\", \"977\u00d75=4885\"),\n    @(\"715\u00d73=2145\", \"599\u00d78=4792\"),\n    @(\"725\u00d76=4350\", \"216\u00d79=1944\"),\n    @(\"692\u00d74=2768\", \"902\u00d76=5412\"),\n    @(\"305\u00d72=610\",  \"824\u00d77=5768\"),\n    @(\"443\u00d79=3987\", \"460\u00d79=4140\"),\n    @(\"420\u00d75=2100\", \"833\u00d73=2499\"),\n    @(\"885\u00d79=7965\", \"311\u00d73=933\"),\n    @(\"984\u00d72=1968\", \"803\u00d78=6424\"),\n    @(\"529\u00d77=3703\", \"912\u00d79=8208\"),\n    @(\"340\u00d72=680\",  \"223\u00d79=2007\"),\n    @(\"309\u00d74=1236\", \"987\u00d73=2961\"),\n    @(\"107\u00d72=214\",  \"318\u00d73=954\"),\n    @(\"436\u00d77=3052\", \"799\u00d74=3196\"),\n    @(\"768\u00d77=5376\", \"293\u00d75=1465\"),\n    @(\"628\u00d74=2512\", \"690\u00d73=2070\"),\n    @(\"852\u00d76=5112\", \"748\u00d76=4488\"),\n    @(\"132\u00d79=1188\", \"956\u00d76=5736\"),\n    @(\"185\u00d75=925\",  \"350\u00d73=1050\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
